# Generate Report for Handoff
# New handoff run produced a new GUID-named source file / hash and updated
# handoff timestamps; update every sheet that mirrors this info.

$wb = $excel.ActiveWorkbook

$oldGuid = "b19c0934-1ab3-4f08-b96e-121239b6953b"
$newGuid = "9ab5c9c0-b015-42f7-a9ac-760b9551fc85"
$oldHash = "058cf82b8649910534183c5ee074b8a189c5d4f0"
$newHash = "da7ebbbd8b48b4da56778c172843ad198defc694"
$commitSha = "6e11991caac89115fb90809e1dde87acbd13ef7c"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-31 15:15:40"

$overviewTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md"
$overviewDisplay = "e2e\$newGuid.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewTarget, "", "", $overviewDisplay)

# --- zh-cn sheet ---
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-31 15:15:35"

$zhTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md"
$zhDisplay = "$newGuid.md"
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhTarget, "", "", $zhDisplay)

# --- de-de sheet ---
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-31 15:15:40"

$deTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md"
$deDisplay = "$newGuid.md"
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deTarget, "", "", $deDisplay)
